# natmiOut/OldD7/LR-pairs_lrc2p/Ccl21b-Ccr7.xlsx
#
# "Natmi following Dr Hou advice": the single-row ligand/receptor edge
# table (Ccl21b -> Ccr7) is re-run with an expanded set of sending /
# target clusters. The former single data row (sCs -> M2) is replaced
# by six rows covering every combination of the two new sending
# clusters (FAPs, sCs) against the three target clusters (ECs, FAPs,
# M2), each carrying its own freshly computed NATMI specificity scores.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ccl21b"
$ws.Range("C2").Value = "Ccr7"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.134289
$ws.Range("H2").Value = 0.402867
$ws.Range("I2").Value = 0.3678949098679525
$ws.Range("J2").Value = 0.3678949098679525
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1277553333333333
$ws.Range("N2").Value = 0.383266
$ws.Range("O2").Value = 0.04516525669351801
$ws.Range("P2").Value = 0.04516525669351801
$ws.Range("Q2").Value = 0.017156135958
$ws.Range("R2").Value = 0.154405223622
$ws.Range("S2").Value = 0.01661606804042474
$ws.Range("T2").Value = 0.01661606804042474

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ccl21b"
$ws.Range("C3").Value = "Ccr7"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.134289
$ws.Range("H3").Value = 0.402867
$ws.Range("I3").Value = 0.3678949098679525
$ws.Range("J3").Value = 0.3678949098679525
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.186073
$ws.Range("N3").Value = 0.558219
$ws.Range("O3").Value = 0.06578226199610435
$ws.Range("P3").Value = 0.06578226199610435
$ws.Range("Q3").Value = 0.02498755709700001
$ws.Range("R3").Value = 0.224888013873
$ws.Range("S3").Value = 0.02420095934796684
$ws.Range("T3").Value = 0.02420095934796684

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ccl21b"
$ws.Range("C4").Value = "Ccr7"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.134289
$ws.Range("H4").Value = 0.402867
$ws.Range("I4").Value = 0.3678949098679525
$ws.Range("J4").Value = 0.3678949098679525
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.514791333333333
$ws.Range("N4").Value = 7.544373999999999
$ws.Range("O4").Value = 0.8890524813103776
$ws.Range("P4").Value = 0.8890524813103776
$ws.Range("Q4").Value = 0.337708813362
$ws.Range("R4").Value = 3.039379320258
$ws.Range("S4").Value = 0.3270778824795609
$ws.Range("T4").Value = 0.3270778824795609

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Ccl21b"
$ws.Range("C5").Value = "Ccr7"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.230731
$ws.Range("H5").Value = 0.6921929999999999
$ws.Range("I5").Value = 0.6321050901320475
$ws.Range("J5").Value = 0.6321050901320475
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1277553333333333
$ws.Range("N5").Value = 0.383266
$ws.Range("O5").Value = 0.04516525669351801
$ws.Range("P5").Value = 0.04516525669351801
$ws.Range("Q5").Value = 0.02947711581533333
$ws.Range("R5").Value = 0.265294042338
$ws.Range("S5").Value = 0.02854918865309326
$ws.Range("T5").Value = 0.02854918865309326

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Ccl21b"
$ws.Range("C6").Value = "Ccr7"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.230731
$ws.Range("H6").Value = 0.6921929999999999
$ws.Range("I6").Value = 0.6321050901320475
$ws.Range("J6").Value = 0.6321050901320475
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.186073
$ws.Range("N6").Value = 0.558219
$ws.Range("O6").Value = 0.06578226199610435
$ws.Range("P6").Value = 0.06578226199610435
$ws.Range("Q6").Value = 0.042932809363
$ws.Range("R6").Value = 0.386395284267
$ws.Range("S6").Value = 0.0415813026481375
$ws.Range("T6").Value = 0.0415813026481375

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Ccl21b"
$ws.Range("C7").Value = "Ccr7"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.230731
$ws.Range("H7").Value = 0.6921929999999999
$ws.Range("I7").Value = 0.6321050901320475
$ws.Range("J7").Value = 0.6321050901320475
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.514791333333333
$ws.Range("N7").Value = 7.544373999999999
$ws.Range("O7").Value = 0.8890524813103776
$ws.Range("P7").Value = 0.8890524813103776
$ws.Range("Q7").Value = 0.5802403191313332
$ws.Range("R7").Value = 5.222162872181999
$ws.Range("S7").Value = 0.5619745988308167
$ws.Range("T7").Value = 0.5619745988308167
